$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: 1602 -> 3061
$ws.Range("B2").Value = 3061

# Update A3: 2 -> 1, B3: 1459 -> 1197 (values from old row 4 shifted up)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1197

# Remove old row 4 entirely (it gets deleted, shrinking the dimension to A1:B3)
$ws.Rows("4").Delete()
